$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Insert a new row above row 5 (shifts existing rows 5..166 -> 6..167)
# ---------------------------------------------------------------------
$ws.Rows.Item(5).Insert()

# Copy the formatting from row 6 (the row that used to be row 5) into the
# newly inserted row 5 for the H/I columns so the styles match.
$ws.Range("H6:I6").Copy()
$ws.Range("H5:I5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New "interest" entry inserted into the H/I side list.
$ws.Range("H5").Value = "13/08/2024"
$ws.Range("I5").Value = 20

# ---------------------------------------------------------------------
# 2) Append two new transaction rows at the bottom of the ledger
#    (rows 168 and 169).
# ---------------------------------------------------------------------
$ws.Range("A166:D167").Copy()
$ws.Range("A168:D169").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 168: reuse the exact text already used elsewhere in the workbook so
# the shared-string table matches (same Vietnamese text, same encoding).
$interestText = $ws.Range("B38").Value2
$ws.Range("A168").Value = "13/08/2024"
$ws.Range("B168").Value = $interestText
$ws.Range("C168").Value = 3000
$ws.Range("D168").Formula = "=D167+C168"

# Row 169: no value in column A.
$loanText = $ws.Range("B19").Value2.Replace("10tr", "20tr")
$ws.Range("A169").Clear()
$ws.Range("B169").Value = $loanText
$ws.Range("C169").Value = -20000
$ws.Range("D169").Formula = "=D168+C169"

# ---------------------------------------------------------------------
# 3) Window/view changes: make "CÔ DIỄM" the active/selected sheet,
#    scroll it to A2, and select M19.
# ---------------------------------------------------------------------
$ws.Activate()
$aw = $excel.ActiveWindow
$aw.ScrollRow = 2
$aw.ScrollColumn = 1
$ws.Range("M19").Select()
